# MERN Project 33 - Menyimpan Postingan Blog Ke Database MongoDB
# The API-planning docs workbook renamed the JSON response field "code" to
# "status" throughout the sample payloads on the "API" sheet, and the
# active/selected sheet + cell changed from "Database"/B1 to "API"/G8.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("API")
$ws4 = $wb.Worksheets.Item("Database")

# --- Rename the "code" JSON key to "status" in every sample request/response
#     body on the API sheet (columns F = response, G = error-response). ---
$ws1.Range("F4").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully registered access"",`n    ""data"":{`n        ""name"":""Testing"",`n        ""email"":""test@gmail.com""`n    }`n}"
$ws1.Range("G4").Value = "{`n    ""status"":""400"",`n    ""message"":""email tidak valid""`n}"
$ws1.Range("F5").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully login"",`n    ""data"":{`n        ""id"":""123"",`n        ""token"":""G3121NJ1338DSN""`n    }`n}"
$ws1.Range("G5").Value = "{`n    ""status"":""400"",`n    ""message"":""Password salah""`n}"
$ws1.Range("F6").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully logout""`n}"
$ws1.Range("G6").Value = "{`n    ""status"":""500"",`n    ""message"":""Internal Server Error""`n}"
$ws1.Range("F8").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully added data"",`n    ""data"": {`n            ""title_post"": ""Harry potter book 1"",  `n            ""body_post"": ""Once upon a time bla bla bla"",`n            ""thumb_image"": ""image.JPG"",`n            ""author"": {`n                ""user_id"": 1,`n                ""name"": ""Testing""`n            }`n    }`n}"
$ws1.Range("G8").Value = "{`n    ""status"":""400"",`n    ""message"":""Judul post sudah ada""`n}"
$ws1.Range("F9").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully updated data""`n    ""data"": {}`n}"
$ws1.Range("G9").Value = "{`n    ""status"":""500"",`n    ""message"":""Internal Server Error""`n}"
$ws1.Range("F10").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully updated data""`n}"
$ws1.Range("G10").Value = "{`n    ""status"":""400"",`n    ""message"":""Body post tidak boleh memuat karakter khusus""`n}"
$ws1.Range("F11").Value = "{`n    ""status"":""200"",`n    ""message"":""Successfully deleted data""`n}"
$ws1.Range("G11").Value = "{`n    ""status"":""400"",`n    ""message"":""Id post tidak ditemukan""`n}"

# --- View state: the "Database" sheet was the active/selected tab with B1
#     selected; the edit moves the active tab + selection to "API" at G8
#     (keeping the API sheet's existing D5 scroll position), leaving
#     "Database" deselected with G13 as its last selection. ---
$ws4.Range("G13").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 4
$ws1.Range("G8").Select()
